$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F updates for rows 3-20
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 135
$ws1.Range("F4").Value = 2093
$ws1.Range("F6").Value = 649
$ws1.Range("F8").Value = 2082
$ws1.Range("F9").Value = 10754
$ws1.Range("F12").Value = 288
$ws1.Range("F13").Value = 204
$ws1.Range("F14").Value = 422
$ws1.Range("F15").Value = 8686
$ws1.Range("F16").Value = 1115
$ws1.Range("F17").Value = 727
$ws1.Range("F18").Value = 3419
$ws1.Range("F19").Value = 70
$ws1.Range("F20").Value = 3349

# Sheet "全部类型" (sheet4) - column F updates for rows 3-23
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 135
$ws4.Range("F4").Value = 2093
$ws4.Range("F6").Value = 649
$ws4.Range("F9").Value = 2082
$ws4.Range("F12").Value = 10754
$ws4.Range("F15").Value = 288
$ws4.Range("F16").Value = 204
$ws4.Range("F17").Value = 422
$ws4.Range("F18").Value = 8686
$ws4.Range("F19").Value = 1115
$ws4.Range("F20").Value = 727
$ws4.Range("F21").Value = 3419
$ws4.Range("F22").Value = 70
$ws4.Range("F23").Value = 3349
